$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 78
$ws.Range("I6").Value = 78
$ws.Range("K6").Value = 234
$ws.Range("M6").Value = -122

$ws.Range("H19").Value = 2179.5
$ws.Range("I19").Value = 4338.5
$ws.Range("K19").Value = 4338.5
$ws.Range("M19").Value = -4163.5

$ws.Range("H33").Value = 333540.47
$ws.Range("I33").Value = 357357.72
$ws.Range("J33").Value = 99
$ws.Range("K33").Value = 357357.72
$ws.Range("L33").Value = 99
$ws.Range("M33").Value = -357128.72
$ws.Range("N33").Value = -557

$ws.Range("H51").Value = 6324.6
$ws.Range("I51").Value = 5408.3335
$ws.Range("J51").Value = 6935.4443
$ws.Range("K51").Value = 5408.3335
$ws.Range("L51").Value = 6935.4443
$ws.Range("M51").Value = -4924.3335
$ws.Range("N51").Value = -7903.4443

$ws.Range("H100").Value = 886.8570999999999
$ws.Range("I100").Value = 1003
$ws.Range("J100").Value = 732
$ws.Range("K100").Value = 1003
$ws.Range("L100").Value = 732
$ws.Range("M100").Value = -462
$ws.Range("N100").Value = -1814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 728033.9399999999
$ws.Range("I2").Value = 1164058.1
$ws.Range("K2").Value = 1164058.1
$ws.Range("M2").Value = -1163945.1

$ws.Range("H45").Value = 4078.1177
$ws.Range("I45").Value = 3885.8333
$ws.Range("K45").Value = 3885.8333
$ws.Range("M45").Value = -3508.8333

$ws.Range("H74").Value = 15627293
$ws.Range("I74").Value = 27779456
$ws.Range("K74").Value = 27779456
$ws.Range("M74").Value = -27778582

$ws.Range("H77").Value = 15627293
$ws.Range("I77").Value = 27779456
$ws.Range("K77").Value = 138897280
$ws.Range("M77").Value = -138892912

$ws.Range("H97").Value = 1784.7142
$ws.Range("I97").Value = 2748.25
$ws.Range("K97").Value = 2748.25
$ws.Range("M97").Value = -2252.25

$ws.Range("H102").Value = 327705.75
$ws.Range("I102").Value = 490373.25
$ws.Range("J102").Value = 2370.7144
$ws.Range("K102").Value = 490373.25
$ws.Range("L102").Value = 2370.7144
$ws.Range("M102").Value = -488751.25
$ws.Range("N102").Value = -5614.7144

$ws.Range("H110").Value = 4078.875
$ws.Range("J110").Value = 6333.6665
$ws.Range("L110").Value = 6333.6665
$ws.Range("N110").Value = -10423.6665

$ws.Range("H116").Value = 728033.9399999999
$ws.Range("I116").Value = 1164058.1
$ws.Range("K116").Value = 1164058.1
$ws.Range("M116").Value = -1161764.1

$ws.Range("H132").Value = 5498.75
$ws.Range("I132").Value = 2448.0625
$ws.Range("K132").Value = 7344.1875
$ws.Range("M132").Value = -4814.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 728033.9399999999
$ws.Range("I3").Value = 1164058.1
$ws.Range("K3").Value = 1164058.1
$ws.Range("M3").Value = -1163944.1

$ws.Range("H20").Value = 3022.1875
$ws.Range("I20").Value = 2311.0715
$ws.Range("K20").Value = 2311.0715
$ws.Range("M20").Value = -2064.0715

$ws.Range("H86").Value = 1735.3636
$ws.Range("J86").Value = 1499.75
$ws.Range("L86").Value = 1499.75
$ws.Range("N86").Value = -3745.75

$ws.Range("H89").Value = 1735.3636
$ws.Range("J89").Value = 1499.75
$ws.Range("L89").Value = 7498.75
$ws.Range("N89").Value = -18730.75

$ws.Range("H94").Value = 623617.5600000001
$ws.Range("I94").Value = 761753.75
$ws.Range("J94").Value = 2004.5
$ws.Range("K94").Value = 761753.75
$ws.Range("L94").Value = 2004.5
$ws.Range("M94").Value = -761302.75
$ws.Range("N94").Value = -2906.5

$ws.Range("H107").Value = 1780
$ws.Range("J107").Value = 2000
$ws.Range("L107").Value = 2000
$ws.Range("N107").Value = -5840

$ws.Range("H134").Value = 6111.375
$ws.Range("I134").Value = 3226.75
$ws.Range("K134").Value = 9680.25
$ws.Range("M134").Value = -7145.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4747.636
$ws.Range("I58").Value = 4768.6665
$ws.Range("J58").Value = 4722.4
$ws.Range("K58").Value = 4768.6665
$ws.Range("L58").Value = 4722.4
$ws.Range("M58").Value = -4565.6665
$ws.Range("N58").Value = -5128.4

$ws.Range("H136").Value = 4747.636
$ws.Range("I136").Value = 4768.6665
$ws.Range("J136").Value = 4722.4
$ws.Range("K136").Value = 14305.9995
$ws.Range("L136").Value = 14167.2
$ws.Range("M136").Value = -11755.9995
$ws.Range("N136").Value = -19267.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1750
$ws.Range("J60").Value = 1750
$ws.Range("L60").Value = 5250
$ws.Range("N60").Value = -5752

$ws.Range("H97").Value = 214.28572
$ws.Range("J97").Value = 199.75
$ws.Range("L97").Value = 599.25
$ws.Range("N97").Value = -1591.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 2700
$ws.Range("J33").Value = 2960
$ws.Range("L33").Value = 2960
$ws.Range("N33").Value = -3464

$ws.Range("H97").Value = 255.6
$ws.Range("J97").Value = 290
$ws.Range("L97").Value = 290
$ws.Range("N97").Value = -1282

$ws.Range("H102").Value = 13894901
$ws.Range("I102").Value = 20005968
$ws.Range("J102").Value = 6108.8184
$ws.Range("K102").Value = 20005968
$ws.Range("L102").Value = 6108.8184
$ws.Range("M102").Value = -20004346
$ws.Range("N102").Value = -9352.8184

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1152.5294
$ws.Range("I22").Value = 999.5
$ws.Range("J22").Value = 1371.1428
$ws.Range("K22").Value = 999.5
$ws.Range("L22").Value = 1371.1428
$ws.Range("M22").Value = -704.5
$ws.Range("N22").Value = -1961.1428

$ws.Range("H27").Value = 1152.5294
$ws.Range("I27").Value = 999.5
$ws.Range("J27").Value = 1371.1428
$ws.Range("K27").Value = 999.5
$ws.Range("L27").Value = 1371.1428
$ws.Range("M27").Value = -892.5
$ws.Range("N27").Value = -1585.1428

$ws.Range("H46").Value = 6063.609
$ws.Range("J46").Value = 6918
$ws.Range("L46").Value = 6918
$ws.Range("N46").Value = -7294

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H122").Value = 4470697.5
$ws.Range("I122").Value = 3464.9048
$ws.Range("K122").Value = 10394.7144
$ws.Range("M122").Value = -7944.714399999999

$ws.Range("H132").Value = 3188.138
$ws.Range("I132").Value = 2418.8
$ws.Range("K132").Value = 7256.400000000001
$ws.Range("M132").Value = -4726.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H107").Value = 3166.087
$ws.Range("I107").Value = 3380.8
$ws.Range("J107").Value = 2763.5
$ws.Range("K107").Value = 10142.4
$ws.Range("L107").Value = 8290.5
$ws.Range("M107").Value = -8222.400000000001
$ws.Range("N107").Value = -12130.5

$ws.Range("H113").Value = 825.5
$ws.Range("J113").Value = 977.7778
$ws.Range("L113").Value = 2933.3334
$ws.Range("N113").Value = -7273.3334

$ws.Range("H116").Value = 79927.836
$ws.Range("J116").Value = 79927.836
$ws.Range("L116").Value = 79927.836
$ws.Range("N116").Value = -89105.836

$ws.Range("H122").Value = 3955.843
$ws.Range("I122").Value = 3394.4722
$ws.Range("K122").Value = 10183.4166
$ws.Range("M122").Value = -7733.4166

$ws.Range("H136").Value = 10000.508
$ws.Range("I136").Value = 5606.9473
$ws.Range("K136").Value = 16820.8419
$ws.Range("M136").Value = -14270.8419
